$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range of data rows (row 19 through the row before the
# trailing summary block). We look at column A text values and replace the
# patch-type suffixes "_Fixed" -> "_Manual" and "_Repaired" -> "_Auto" both
# in the literal cell text (column A) and in the SEARCH() formulas (column J)
# that classify each row by its suffix.

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string]) {
        if ($val -like "*_Fixed") {
            $cell.Value2 = $val.Substring(0, $val.Length - 6) + "_Manual"
        } elseif ($val -like "*_Repaired") {
            $cell.Value2 = $val.Substring(0, $val.Length - 9) + "_Auto"
        }
    }

    $jcell = $ws.Cells.Item($r, 10)
    $formula = $jcell.Formula
    if ($formula -ne $null -and $formula -is [string] -and $formula.StartsWith("=") -and $formula.IndexOf("_Fixed") -ge 0) {
        $newFormula = $formula.Replace("*_Fixed", "*_Manual").Replace("""Fixed""", """Manual""").Replace("*_Repaired", "*_Auto").Replace("""Repaired""", """Auto""")
        $jcell.Formula = $newFormula
    }
}

$ws.Range("E7").Select()
